# Rename the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Dist_Widgets"

# Apply 0.00 number format to C4:E4
$ws.Range("C4:E4").NumberFormat = "0.00"

# Add new cell E8 with same number format, empty value, selected
$ws.Range("E8").NumberFormat = "0.00"
$ws.Range("E8").Select()
